# Refresh the cryptos list (GitHub Actions scrape) with updated prices / 1h
# volume deltas, plus a couple of ranking swaps (XRP<->BNB, ImmutableX<->
# LidoDAOToken, BabyDogeCoin<->RocketPoolETH).
#
# NOTE: the Price column holds plain text in the source workbook (e.g.
# "1.000", "0.7015") so values that would otherwise parse as a clean number
# are written with a leading apostrophe -- the normal Excel "force text"
# idiom -- to keep them as literal strings instead of being coerced to
# numeric 1, 0.7015, etc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.195.19'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '1.859.88'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").Value = '''0.7015'
$ws.Range("E5").Value = '  -2.16%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '''241.67'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '''0.07807'
$ws.Range("E8").Value = '  -2.11%  '
$ws.Range("D9").Value = '''0.3105'
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").Value = '''23.82'
$ws.Range("E10").Value = '  -4.34%  '
$ws.Range("D11").Value = '''0.07797'
$ws.Range("E11").Value = '  -3.51%  '
$ws.Range("D12").Value = '1.842.41'
$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("D13").Value = '''92.55'
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("D14").Value = '''5.115'
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("D15").Value = '''0.6888'
$ws.Range("E15").Value = '  -2.72%  '
$ws.Range("D16").Value = '''6.526'
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").Value = '''0.000008436'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '29.190.90'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = '''249.82'
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").Value = '2.107.14'
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").Value = '''12.90'
$ws.Range("E21").Value = '  -3.25%  '
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '''7.582'
$ws.Range("E23").Value = '  -1.29%  '
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = '''0.1530'
$ws.Range("E25").Value = '  -3.13%  '
$ws.Range("D26").Value = '''160.30'
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("D27").Value = '''8.875'
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("D28").Value = '''18.56'
$ws.Range("E28").Value = '  -2.17%  '
$ws.Range("D29").Value = '''1.565'
$ws.Range("E29").Value = '  +3.59%  '
$ws.Range("D30").Value = '''4.273'
$ws.Range("E30").Value = '  -3.23%  '
$ws.Range("E31").Value = '  -1.59%  '
$ws.Range("D32").Value = '''1.209'
$ws.Range("E32").Value = '  -1.81%  '
$ws.Range("D33").Value = '''0.05217'
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7574'
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '''1.873'
$ws.Range("E35").Value = '  -3.45%  '
$ws.Range("D36").Value = '''1.175'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = '''2.715'
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("D38").Value = '''0.01859'
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("D39").Value = '1.224.28'
$ws.Range("E39").Value = '  -5.15%  '
$ws.Range("D40").Value = '''2.722'
$ws.Range("E40").Value = '  -1.61%  '
$ws.Range("D41").Value = '''0.9013'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = '''109.61'
$ws.Range("E42").Value = '  -1.69%  '
$ws.Range("D43").Value = '''5.764'
$ws.Range("E43").Value = '  -10.01%  '
$ws.Range("D44").Value = '''1.000'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '''0.00000000124'
$ws.Range("E45").Value = '  -4.03%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.996.37'
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("D47").Value = '''65.25'
$ws.Range("E47").Value = '  -11.94%  '
$ws.Range("D48").Value = '''0.5185'
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").Value = '''9.515'
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").Value = '''1.765'
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("D51").Value = '''7.028'
$ws.Range("E51").Value = '  -0.98%  '
